# ---------------------------------------------------------------------------
# Applies the "Adding licence permits ..." commit to
# apiary_authority_permit_template_v2.docx
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Helper: split the text that currently occupies [$startPos, $startPos + total
# length of $parts) into one run per entry of $parts, by toggling a
# (no-visual-effect) character formatting property on exactly each sub-range.
# Word's COM layer always creates a fresh run at a property-change boundary,
# which is how real Word Find&Replace / formatting ends up splitting runs.
function Split-Runs($doc, $startPos, $parts) {
    $pos = $startPos
    foreach ($part in $parts) {
        $len = $part.Length
        $endPos = $pos + $len
        $rng = $doc.Range($pos, $endPos)
        $orig = $rng.Bold
        $rng.Bold = 1
        $rng.Bold = $orig
        $pos = $pos + $len
    }
}

# ---------------------------------------------------------------------------
# 1. Top line-shape ("Line 21") nudges slightly: position/extent tweak.
# ---------------------------------------------------------------------------
$shape = $d.Shapes.Item(1)
$shape.Top = 2.7
$shape.Left = 6
$shape.Width = 468.55
$shape.Height = 0.6

# ---------------------------------------------------------------------------
# 2. "SCHEDULE 1" + " " (two runs) -> single run "SCHEDULE 1 "
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("SCHEDULE 1 ", $true, $false, $false, $false, `
    $false, $true, 1, $false, "SCHEDULE 1 ", 2)

# ---------------------------------------------------------------------------
# 3. "Apiary Sites" (schedule-1 heading only) -> "Apiary Licensed Sites"
#    Scope the Find to the one paragraph whose whole text is "Apiary Sites"
#    so "Apiary Sites Conditions" elsewhere is left untouched.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Apiary Sites`r") {
        $null = $para.Range.Find.Execute("Apiary Sites", $true, $false, `
            $false, $false, $false, $true, 1, $false, "Apiary Licensed Sites", 2)
    }
}

# ---------------------------------------------------------------------------
# 4. Table grid / cell widths: 1826 -> 1825 and 1573 -> 1574 (twips) on the
#    Apiary Sites schedule table (first table in the document).
# ---------------------------------------------------------------------------
$tbl = $d.Tables.Item(1)
$tbl.Columns.Item(1).Width = 91.25
$tbl.Columns.Item(2).Width = 78.7

# ---------------------------------------------------------------------------
# 5. Swap the loop variable names:
#    {%tr for site in apiary_sites %}          -> apiary_licensed_sites
#    {% for site in apiary_licensed_sites %}    -> apiary_sites
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("{%tr for site in apiary_sites %} ", $true, `
    $false, $false, $false, $false, $true, 1, $false, `
    "{%tr for site in apiary_licensed_sites %} ", 2)

$null = $d.Content.Find.Execute("{% for site in apiary_licensed_sites %}  ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "{% for site in apiary_sites %}  ", 2)

# ---------------------------------------------------------------------------
# 6. SITE DETAILS block: replace the placeholder angle-bracket tokens with
#    Jinja expressions, broken out into several runs (mirrors what Word
#    produces after a sequence of in-place edits/formatting touches).
# ---------------------------------------------------------------------------

# -- <map ref> ---------------------------------------------------------------
$find = $d.Content
$null = $find.Find.Execute("<map ref>", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
$anchorStart = $find.Start
$anchorEnd = $find.End

$mapStart = $anchorStart - 1          # include the tab right before
$mapEnd = $anchorEnd + 2              # include the two tabs right after

$mapParts = @("`t", "{{ ", "map_ref ", "}}", "`t`t")
$mapFull = [string]::Join("", $mapParts)

$mapRange = $d.Range($mapStart, $mapEnd)
$mapRange.Text = $mapFull
Split-Runs $d $mapStart $mapParts

# -- <forest block> / <COG> / <road/track> -----------------------------------
$find1 = $d.Content
$null = $find1.Find.Execute("<forest block>", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
$fbStart = $find1.Start

$find2 = $d.Content
$null = $find2.Find.Execute("<road/track>", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
$fbEnd = $find2.End

$fbParts = @(
    "{{ ",
    "forest_block ",
    "}}",
    "`t",
    "{{ cog",
    " ",
    "}}",
    "`t`t`t",
    "{{ ",
    "roadtrack ",
    "}}"
)
$fbFull = [string]::Join("", $fbParts)

$fbRange = $d.Range($fbStart, $fbEnd)
$fbRange.Text = $fbFull
Split-Runs $d $fbStart $fbParts

# -- <catchment> / <DRA> ------------------------------------------------------
$find3 = $d.Content
$null = $find3.Find.Execute("<catchment>", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
$catStart = $find3.Start - 1          # include the tab right before

$find4 = $d.Content
$null = $find4.Find.Execute("<DRA>", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
$catEnd = $find4.End + 1              # consume the trailing tab too

$catParts = @(
    "`t",
    "{{ ",
    "catchment ",
    "}}",
    "`t`t",
    "{{ dra_permit }}"
)
$catFull = [string]::Join("", $catParts)

$catRange = $d.Range($catStart, $catEnd)
$catRange.Text = $catFull
Split-Runs $d $catStart $catParts
